# Set default object for all job and level
# Inserts a new "DefaultObject" row at row 11 of the Property1 sheet,
# duplicating the existing row 11 (NPC001) contents/formatting and
# pushing the rest of the table down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property1")

# Duplicate row 11 (NPC001) by copying it and inserting the copy above
# itself; this shifts rows 11..20 down to 12..21 and leaves a clone of
# the original row 11 behind as the new row 11.
$ws.Rows.Item(11).Copy()
$ws.Rows.Item(11).Insert()
$excel.CutCopyMode = $false

# Re-label the new row as the default object entry.
$ws.Range("A11").Value = "DefaultObject"
$ws.Range("J11").Value = "Prefabs/Object/DefaultObject"

# Restore the view: scroll back to column A and select J12 (first prefab
# cell of the row that used to be row 11).
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("J12").Select()
